# Updates cryptos price/volume data, and swaps VeChain/TrustWalletToken row order
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($address, $text) {
    $range = $ws.Range($address)
    # Force text storage so numeric-looking strings (prices, percentages)
    # are not silently reinterpreted by Excel as numbers/dates.
    $range.NumberFormat = "@"
    $range.Value = $text
    # Restore the default (unstyled) cell style so no stray formatting is introduced.
    $range.Style = "Normal"
}

Set-CellText "D2" "27.382.53"
Set-CellText "E2" "  -2.87%  "
Set-CellText "D3" "1.858.41"
Set-CellText "E3" "  -3.73%  "
Set-CellText "E4" "  -0.86%  "
Set-CellText "D5" "323.67"
Set-CellText "E5" "  +0.60%  "
Set-CellText "E6" "  -0.89%  "
Set-CellText "D7" "0.4529"
Set-CellText "E7" "  -4.30%  "
Set-CellText "D8" "0.3873"
Set-CellText "E8" "  -4.39%  "
Set-CellText "D9" "48.91"
Set-CellText "E9" "  -8.39%  "
Set-CellText "D10" "0.07920"
Set-CellText "E10" "  -6.98%  "
Set-CellText "D11" "1.016"
Set-CellText "E11" "  -3.23%  "
Set-CellText "D12" "21.39"
Set-CellText "E12" "  -3.74%  "
Set-CellText "D13" "1.855.01"
Set-CellText "E13" "  -6.22%  "
Set-CellText "D14" "5.921"
Set-CellText "E14" "  -3.17%  "
Set-CellText "D15" "7.133"
Set-CellText "E15" "  -5.12%  "
Set-CellText "D16" "1.004"
Set-CellText "E16" "  -0.94%  "
Set-CellText "D17" "86.04"
Set-CellText "E17" "  -4.25%  "
Set-CellText "D18" "0.00001034"
Set-CellText "E18" "  -3.37%  "
Set-CellText "D19" "0.06518"
Set-CellText "E19" "  -1.38%  "
Set-CellText "E20" "  -6.26%  "
Set-CellText "E21" "  -0.97%  "
Set-CellText "D22" "5.543"
Set-CellText "E22" "  -4.09%  "
Set-CellText "D23" "27.375.60"
Set-CellText "E23" "  -3.07%  "
Set-CellText "D24" "10.86"
Set-CellText "E24" "  -4.95%  "
Set-CellText "D25" "2.281"
Set-CellText "E25" "  -1.02%  "
Set-CellText "D26" "2.076.96"
Set-CellText "E26" "  -5.83%  "
Set-CellText "D27" "153.94"
Set-CellText "E27" "  -0.70%  "
Set-CellText "D28" "19.93"
Set-CellText "E28" "  -1.19%  "
Set-CellText "D29" "2.079"
Set-CellText "E29" "  -3.91%  "
Set-CellText "D30" "5.443"
Set-CellText "E30" "  -5.55%  "
Set-CellText "D31" "121.13"
Set-CellText "E31" "  -2.16%  "
Set-CellText "D32" "1.484"
Set-CellText "E32" "  +2.64%  "
Set-CellText "E33" "  -3.12%  "
Set-CellText "D34" "0.9376"
Set-CellText "E34" "  -4.34%  "
Set-CellText "E35" "  -1.90%  "
Set-CellText "D36" "5.267"
Set-CellText "E36" "  -5.73%  "
Set-CellText "B37" "TrustWalletToken"
Set-CellText "C37" "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-CellText "D37" "1.233"
Set-CellText "E37" "  -0.36%  "
Set-CellText "B38" "VeChain"
Set-CellText "C38" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-CellText "D38" "0.02237"
Set-CellText "E38" "  -3.66%  "
Set-CellText "D39" "0.06001"
Set-CellText "E39" "  -2.79%  "
Set-CellText "D40" "8.199"
Set-CellText "E40" "  -11.66%  "
Set-CellText "E41" "  -0.87%  "
Set-CellText "D42" "0.5917"
Set-CellText "E42" "  -4.38%  "
Set-CellText "D43" "0.1905"
Set-CellText "D44" "10.12"
Set-CellText "E44" "  -9.02%  "
Set-CellText "D45" "1.283"
Set-CellText "E45" "  -2.86%  "
Set-CellText "D46" "0.5622"
Set-CellText "E46" "  -4.79%  "
Set-CellText "D47" "12.01"
Set-CellText "E47" "  -6.68%  "
Set-CellText "D48" "3.373"
Set-CellText "E48" "  -0.66%  "
Set-CellText "E49" "  -5.76%  "
Set-CellText "D50" "0.06770"
Set-CellText "E50" "  +0.02%  "
Set-CellText "D51" "108.48"
Set-CellText "E51" "  -1.27%  "
